# Update countries & provincias Spain
# Applies the data refresh captured in the commit:
#  - Belice and Nueva Caledonia swap their list position (and thus each
#    keeps its own "Casos activos"/"Muertes" figures in rows 192/193)
#  - Several country rows get refreshed case counts
#  - The "last updated" timestamp moves from 19:34 to 20:04

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Swap "Belice" and "Nueva Caledonia" in the country list (rows 192/193) ---
# Row 192 keeps pointing at the (now) "Belice" label; row 193 now holds
# "Nueva Caledonia". The country-name cells swap, and so do the data values
# that travel with each country (only D and H changed between them; the
# other columns already matched).
$ws.Range("A192").Value = "Belice"
$ws.Range("A193").Value = "Nueva Caledonia"

$ws.Range("D192").Value = 16
$ws.Range("H192").Value = 2

$ws.Range("D193").Value = 18
$ws.Range("H193").Value = 0

# --- Refresh the country statistics ---

# Row 4: Estados Unidos
$ws.Range("B4").Value = 1333374
$ws.Range("C4").Value = 11589
$ws.Range("E4").Value = 1029497
$ws.Range("G4").Value = 629
$ws.Range("H4").Value = 79244

# Row 10: Alemania
$ws.Range("B10").Value = 171021
$ws.Range("C10").Value = 433
$ws.Range("E10").Value = 20196
$ws.Range("G10").Value = 15
$ws.Range("H10").Value = 7525

# Row 16: India
$ws.Range("B16").Value = 62769
$ws.Range("C16").Value = 3074
$ws.Range("E16").Value = 42008

# Row 24: Pakistan
$ws.Range("B24").Value = 28736
$ws.Range("C24").Value = 2301
$ws.Range("D24").Value = 7809
$ws.Range("E24").Value = 20291
$ws.Range("G24").Value = 37
$ws.Range("H24").Value = 636

# Row 32: Emiratos Arabes Unidos
$ws.Range("B32").Value = 17417
$ws.Range("C32").Value = 624
$ws.Range("D32").Value = 4295
$ws.Range("E32").Value = 12937
$ws.Range("G32").Value = 11
$ws.Range("H32").Value = 185

# Row 33: Israel
$ws.Range("B33").Value = 16454
$ws.Range("C33").Value = 18
$ws.Range("D33").Value = 11376
$ws.Range("E33").Value = 4831
$ws.Range("G33").Value = 2
$ws.Range("H33").Value = 247

# Row 37: Rumania
$ws.Range("E37").Value = 7283
$ws.Range("G37").Value = 13
$ws.Range("H37").Value = 936

# --- Update the "last updated" timestamp in the title cell (A1) ---
$ws.Range("A1").Value = "Datos actualizados a 9 de Mayo de 2020 a las 20:04"
